# Update Row 3 ("R" - conference round) values on both the OFF and DEF
# sheets to reflect the simulated season logged through the 2021
# divisional round.

$wb = $excel.ActiveWorkbook

# OFF sheet
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 581
$wsOff.Range("C3").Value = 414
$wsOff.Range("D3").Value = 136
$wsOff.Range("E3").Value = 77

# DEF sheet
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 569
$wsDef.Range("C3").Value = 368
$wsDef.Range("D3").Value = 124
$wsDef.Range("E3").Value = 61
